$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency price/volume data (and the MXToken/Aave row swap)

$ws.Range("D2").Value = "27.158.84"
$ws.Range("E2").Value = "  -1.97%  "
$ws.Range("D3").Value = "1.557.90"
$ws.Range("E3").Value = "  -2.17%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "206.42"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.80%  "
$ws.Range("E6").Value = "  -2.51%  "
$ws.Range("E8").Value = "  -0.63%  "
$ws.Range("E9").Value = "  -2.15%  "
$ws.Range("E10").Value = "  -0.10%  "
$ws.Range("E11").Value = "  -0.86%  "
$ws.Range("D12").Value = "1.779.77"
$ws.Range("E12").Value = "  -2.08%  "
$ws.Range("D13").Value = "1.553.09"
$ws.Range("E13").Value = "  -2.46%  "
$ws.Range("E14").Value = "  -2.53%  "
$ws.Range("E15").Value = "  -3.01%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "62.84"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -0.97%  "
$ws.Range("D17").Value = "27.155.88"
$ws.Range("E17").Value = "  -1.90%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "216.15"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -1.84%  "
$ws.Range("D19").Value = "0.0₃0686"
$ws.Range("E19").Value = "  -1.54%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.23"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -1.63%  "
$ws.Range("E21").Value = "  +0.07%  "
$ws.Range("E22").Value = "  -1.17%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.34"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -3.39%  "
$ws.Range("E24").Value = "  +0.13%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "151.99"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -1.29%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.57"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -3.19%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "14.91"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -1.56%  "
$ws.Range("E28").Value = "  -0.01%  "
$ws.Range("E30").Value = "  -1.55%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0462"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -1.87%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.17"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -1.92%  "
$ws.Range("D33").Value = "1.378.56"
$ws.Range("E33").Value = "  +0.31%  "
$ws.Range("E34").Value = "  -0.81%  "
$ws.Range("E35").Value = "  -0.05%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.947"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -3.00%  "
$ws.Range("E37").Value = "  -1.83%  "
$ws.Range("E38").Value = "  -1.47%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.811"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -2.01%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.514"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -4.39%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.985"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +1.80%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.80"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +3.95%  "
$ws.Range("B44").Value = "Aave"
$ws.Range("C44").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "63.25"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -2.14%  "
$ws.Range("B45").Value = "MXToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.16"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -0.31%  "
$ws.Range("E46").Value = "  +0.03%  "
$ws.Range("D47").Value = "1.692.52"
$ws.Range("E47").Value = "  -2.04%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "85.41"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -2.03%  "
$ws.Range("D49").Value = "0.0₇0981"
$ws.Range("E49").Value = "  -2.92%  "
$ws.Range("E50").Value = "  -0.73%  "
$ws.Range("E51").Value = "  +0.09%  "
